$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "gemini-1.5-pro"
$ws.Range("C9").Value = "0.04 ± 0.83"
$ws.Range("D9").Value = "0.26 ± 0.61"
$ws.Range("J9").Value = "0.64 ± 0.37"
$ws.Range("K9").Value = "0.67 ± 0.38"
$ws.Range("L9").Value = "0.66 ± 0.38"
$ws.Range("M9").Value = "0.66 ± 0.38"
$ws.Range("N9").Value = "0.75 ± 0.43"
$ws.Range("P9").Value = "0.45 ± 0.29"
$ws.Range("Q9").Value = "2.45 ± 1.12"
$ws.Range("R9").Value = "0.248 ± 0.00"
$ws.Range("S9").Value = "0.7 ± 0.41"
$ws.Range("T9").Value = "0.75 ± 0.43"
$ws.Range("U9").Value = "2.03 ± 1.52"
$ws.Range("V9").Value = "0.65 ± 0.45"
$ws.Range("W9").Value = "0.7 ± 0.41"
$ws.Range("X9").Value = "0.96 ± 0.59"
